{"js": "// Update the date line and the 25 multiplication answers to the new\n// values from the commit's generated output.\nconst replacements = [\n  [\"2024-06-17 Monday\", \"2024-06-18 Tuesday\"],\n  [\"28\u00d748=1344\", \"77\u00d722=1694\"],\n  [\"85\u00d759=5015\", \"86\u00d720=1720\"],\n  [\"49\u00d761=2989\", \"42\u00d777=3234\"],\n  [\"72\u00d715=1080\", \"27\u00d752=1404\"],\n  [\"46\u00d746=2116\", \"13\u00d784=1092\"],\n  [\"29\u00d765=1885\", \"22\u00d745=990\"],\n  [\"15\u00d780=1200\", \"68\u00d737=2516\"],\n  [\"69\u00d770=4830\", \"57\u00d746=2622\"],\n  [\"28\u00d753=1484\", \"95\u00d744=4180\"],\n  [\"51\u00d731=1581\", \"84\u00d743=3612\"],\n  [\"66\u00d766=4356\", \"24\u00d775=1800\"],\n  [\"66\u00d769=4554\", \"55\u00d769=3795\"],\n  [\"34\u00d752=1768\", \"59\u00d773=4307\"],\n  [\"18\u00d746=828\", \"94\u00d790=8460\"],\n  [\"42\u00d778=3276\", \"13\u00d720=260\"],\n  [\"71\u00d752=3692\", \"68\u00d750=3400\"],\n  [\"18\u00d790=1620\", \"37\u00d714=518\"],\n  [\"22\u00d738=836\", \"75\u00d732=2400\"],\n  [\"68\u00d775=5100\", \"39\u00d725=975\"],\n  [\"64\u00d750=3200\", \"44\u00d784=3696\"],\n  [\"73\u00d738=2774\", \"65\u00d781=5265\"],\n  [\"66\u00d745=2970\", \"50\u00d750=2500\"],\n  [\"14\u00d771=994\", \"80\u00d714=1120\"],\n  [\"18\u00d773=1314\", \"47\u00d763=2961\"],\n  [\"93\u00d741=3813\", \"32\u00d721=672\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and the 25 multiplication answers to the new\n# values from the commit's generated output.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-06-17 Monday\", \"2024-06-18 Tuesday\"),\n    @(\"28\u00d748=1344\", \"77\u00d722=1694\"),\n    @(\"85\u00d759=5015\", \"86\u00d720=1720\"),\n    @(\"49\u00d761=2989\", \"42\u00d777=3234\"),\n    @(\"72\u00d715=1080\", \"27\u00d752=1404\"),\n    @(\"46\u00d746=2116\", \"13\u00d784=1092\"),\n    @(\"29\u00d765=1885\", \"22\u00d745=990\"),\n    @(\"15\u00d780=1200\", \"68\u00d737=2516\"),\n    @(\"69\u00d770=4830\", \"57\u00d746=2622\"),\n    @(\"28\u00d753=1484\", \"95\u00d744=4180\"),\n    @(\"51\u00d731=1581\", \"84\u00d743=3612\"),\n    @(\"66\u00d766=4356\", \"24\u00d775=1800\"),\n    @(\"66\u00d769=4554\", \"55\u00d769=3795\"),\n    @(\"34\u00d752=1768\", \"59\u00d773=4307\"),\n    @(\"18\u00d746=828\", \"94\u00d790=8460\"),\n    @(\"42\u00d778=3276\", \"13\u00d720=260\"),\n    @(\"71\u00d752=3692\", \"68\u00d750=3400\"),\n    @(\"18\u00d790=1620\", \"37\u00d714=518\"),\n    @(\"22\u00d738=836\", \"75\u00d732=2400\"),\n    @(\"68\u00d775=5100\", \"39\u00d725=975\"),\n    @(\"64\u00d750=3200\", \"44\u00d784=3696\"),\n    @(\"73\u00d738=2774\", \"65\u00d781=5265\"),\n    @(\"66\u00d745=2970\", \"50\u00d750=2500\"),\n    @(\"14\u00d771=994\", \"80\u00d714=1120\"),\n    @(\"18\u00d773=1314\", \"47\u00d763=2961\"),\n    @(\"93\u00d741=3813\", \"32\u00d721=672\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $range = $d.Content\n    $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
